$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 339
$ws.Range("F7").Value = 885
$ws.Range("F8").Value = 62
$ws.Range("F9").Value = 530
$ws.Range("F12").Value = 1162
$ws.Range("F15").Value = 41
$ws.Range("F17").Value = 6703
$ws.Range("F20").Value = 21
$ws.Range("F21").Value = 7613
$ws.Range("F24").Value = 3413
$ws.Range("F26").Value = 2136
$ws.Range("F27").Value = 910
$ws.Range("F29").Value = 164
$ws.Range("F35").Value = 1740
$ws.Range("F37").Value = 187
$ws.Range("F38").Value = 58
$ws.Range("F41").Value = 1236
$ws.Range("F42").Value = 1833

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 53

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1237

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1237
$ws.Range("F7").Value = 339
$ws.Range("F9").Value = 885
$ws.Range("F10").Value = 62
$ws.Range("F11").Value = 530
$ws.Range("F14").Value = 1162
$ws.Range("F18").Value = 41
$ws.Range("F20").Value = 6703
$ws.Range("F23").Value = 21
$ws.Range("F24").Value = 7613
$ws.Range("F27").Value = 3413
$ws.Range("F29").Value = 2136
$ws.Range("F30").Value = 910
$ws.Range("F32").Value = 164
$ws.Range("F35").Value = 53
$ws.Range("F38").Value = 1740
$ws.Range("F40").Value = 187
$ws.Range("F41").Value = 58
$ws.Range("F44").Value = 1236
$ws.Range("F45").Value = 1833
